$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O:P. This shifts the existing O..U columns to Q..W,
# carrying their formatting/values/styles along automatically.
$ws.Columns("O:P").Insert()

# Rename the (now shifted-in-place) headers at M1/N1.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Populate the two brand-new header cells, matching the header style (s="1")
# used by every other header cell in row 1 (copy formatting from N1, the
# neighboring original header cell).
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill in the data for the two new columns (O = Correct Pred Predicates
# Parents, P = Correct Pred Predicates Related) per row.
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 3

$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 3

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
